# Auto-generated Excel COM-interop script applying the Shinryu_Profits.xlsx diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 64: H64, J64, L64, N64
$ws.Range("H64").Value = 3916.0193
$ws.Range("J64").Value = 4011.1
$ws.Range("L64").Value = 4011.1
$ws.Range("N64").Value = -4507.1

# Row 67: H67, J67, L67, N67
$ws.Range("H67").Value = 3916.0193
$ws.Range("J67").Value = 4011.1
$ws.Range("L67").Value = 4011.1
$ws.Range("N67").Value = -5727.1

# Row 107: H107, I107, J107, K107, L107, M107, N107
$ws.Range("H107").Value = 5513.522
$ws.Range("I107").Value = 1000.2632
$ws.Range("J107").Value = 26951.5
$ws.Range("K107").Value = 1000.2632
$ws.Range("L107").Value = 26951.5
$ws.Range("M107").Value = 919.7368
$ws.Range("N107").Value = -30791.5

# Row 129: H129
$ws.Range("H129").Value = 707.8421

# Row 138: H138, I138, K138, M138
$ws.Range("H138").Value = 2195.802
$ws.Range("I138").Value = 965
$ws.Range("K138").Value = 2895
$ws.Range("M138").Value = 2245

$ws = $wb.Worksheets.Item("ARM")
# Row 23: I23, J23, K23, L23, M23, N23
$ws.Range("I23").Value = 11248.75
$ws.Range("J23").Value = 10280
$ws.Range("K23").Value = 11248.75
$ws.Range("L23").Value = 10280
$ws.Range("M23").Value = -10989.75
$ws.Range("N23").Value = -10798

# Row 132: H132, I132, J132, K132, L132, M132, N132
$ws.Range("H132").Value = 1609.925
$ws.Range("I132").Value = 926.52
$ws.Range("J132").Value = 2748.9333
$ws.Range("K132").Value = 2779.56
$ws.Range("L132").Value = 8246.7999
$ws.Range("M132").Value = -249.5599999999999
$ws.Range("N132").Value = -13306.7999

$ws = $wb.Worksheets.Item("BSM")
# Row 105: H105, I105, J105, K105, L105, M105, N105
$ws.Range("H105").Value = 1962.5714
$ws.Range("I105").Value = 1411.9
$ws.Range("J105").Value = 2268.5
$ws.Range("K105").Value = 1411.9
$ws.Range("L105").Value = 2268.5
$ws.Range("M105").Value = 335.0999999999999
$ws.Range("N105").Value = -5762.5

# Row 107: H107, I107, K107, M107
$ws.Range("H107").Value = 1508.1818
$ws.Range("I107").Value = 1398.75
$ws.Range("K107").Value = 1398.75
$ws.Range("M107").Value = 521.25

# Row 134: H134, I134, J134, K134, L134, M134, N134
$ws.Range("H134").Value = 2011.6
$ws.Range("I134").Value = 1671.6364
$ws.Range("J134").Value = 4504.6665
$ws.Range("K134").Value = 5014.9092
$ws.Range("L134").Value = 13513.9995
$ws.Range("M134").Value = -2479.9092
$ws.Range("N134").Value = -18583.9995

$ws = $wb.Worksheets.Item("CRP")
# Row 20: H20, J20, L20, N20
$ws.Range("H20").Value = 24689.75
$ws.Range("J20").Value = 24689.75
$ws.Range("L20").Value = 24689.75
$ws.Range("N20").Value = -25161.75

# Row 30: H30, J30, L30, N30
$ws.Range("H30").Value = 24689.75
$ws.Range("J30").Value = 24689.75
$ws.Range("L30").Value = 24689.75
$ws.Range("N30").Value = -24871.75

# Row 31: H31, I31, J31, K31, L31, M31, N31
$ws.Range("H31").Value = 2378.389
$ws.Range("I31").Value = 1395.9706
$ws.Range("J31").Value = 4048.5
$ws.Range("K31").Value = 1395.9706
$ws.Range("L31").Value = 4048.5
$ws.Range("M31").Value = -1100.9706
$ws.Range("N31").Value = -4638.5

# Row 33: H33, I33, J33, K33, L33, M33, N33
$ws.Range("H33").Value = 26147.857
$ws.Range("I33").Value = 2500
$ws.Range("J33").Value = 35607
$ws.Range("K33").Value = 2500
$ws.Range("L33").Value = 35607
$ws.Range("M33").Value = -2121
$ws.Range("N33").Value = -36365

# Row 34: H34, I34, J34, K34, L34, M34, N34
$ws.Range("H34").Value = 2378.389
$ws.Range("I34").Value = 1395.9706
$ws.Range("J34").Value = 4048.5
$ws.Range("K34").Value = 1395.9706
$ws.Range("L34").Value = 4048.5
$ws.Range("M34").Value = -1193.9706
$ws.Range("N34").Value = -4452.5

# Row 58: H58, I58, J58, K58, L58, M58, N58
$ws.Range("H58").Value = 1799.6459
$ws.Range("I58").Value = 1371.2307
$ws.Range("J58").Value = 2305.9546
$ws.Range("K58").Value = 1371.2307
$ws.Range("L58").Value = 2305.9546
$ws.Range("M58").Value = -1168.2307
$ws.Range("N58").Value = -2711.9546

# Row 99: H99, I99, J99, K99, L99, M99, N99
$ws.Range("H99").Value = 2676.5
$ws.Range("I99").Value = 2773.1428
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 2773.1428
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = -1275.1428
$ws.Range("N99").Value = -4996

# Row 123: H123, J123, L123, N123
$ws.Range("H123").Value = 44790
$ws.Range("J123").Value = 44790
$ws.Range("L123").Value = 44790
$ws.Range("N123").Value = -54590

# Row 125: H125, J125, L125, N125
$ws.Range("H125").Value = 23656.5
$ws.Range("J125").Value = 23656.5
$ws.Range("L125").Value = 23656.5
$ws.Range("N125").Value = -28576.5

# Row 126: H126, I126, J126, K126, L126, M126, N126
$ws.Range("H126").Value = 2676.5
$ws.Range("I126").Value = 2773.1428
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 8319.428400000001
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -5849.428400000001
$ws.Range("N126").Value = -10940

# Row 127: H127, J127, L127, N127
$ws.Range("H127").Value = 20975
$ws.Range("J127").Value = 20975
$ws.Range("L127").Value = 20975
$ws.Range("N127").Value = -30895

# Row 128: H128, J128, L128, N128
$ws.Range("H128").Value = 24689.75
$ws.Range("J128").Value = 24689.75
$ws.Range("L128").Value = 24689.75
$ws.Range("N128").Value = -34649.75

# Row 129: H129, J129, L129, N129
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()

# Row 130: H130, J130, L130, N130
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

# Row 131: H131, J131, L131, N131
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

# Row 132: H132, I132, J132, K132, L132, M132, N132
$ws.Range("H132").Value = 1340.8372
$ws.Range("I132").Value = 990.37836
$ws.Range("J132").Value = 3502
$ws.Range("K132").Value = 2971.13508
$ws.Range("L132").Value = 10506
$ws.Range("M132").Value = -441.13508
$ws.Range("N132").Value = -15566

# Row 133: H133, J133, L133, N133
$ws.Range("H133").Value = 20325.5
$ws.Range("J133").Value = 20325.5
$ws.Range("L133").Value = 20325.5
$ws.Range("N133").Value = -25385.5

# Row 134: H134, I134, J134, K134, L134, M134, N134
$ws.Range("H134").Value = 2568.9512
$ws.Range("I134").Value = 1533.3846
$ws.Range("J134").Value = 4363.933
$ws.Range("K134").Value = 4600.1538
$ws.Range("L134").Value = 13091.799
$ws.Range("M134").Value = -2065.1538
$ws.Range("N134").Value = -18161.799

# Row 135: H135, J135, L135, N135
$ws.Range("H135").Value = 36790
$ws.Range("J135").Value = 36790
$ws.Range("L135").Value = 36790
$ws.Range("N135").Value = -46930

# Row 136: H136, I136, J136, K136, L136, M136, N136
$ws.Range("H136").Value = 1799.6459
$ws.Range("I136").Value = 1371.2307
$ws.Range("J136").Value = 2305.9546
$ws.Range("K136").Value = 4113.6921
$ws.Range("L136").Value = 6917.8638
$ws.Range("M136").Value = -1563.6921
$ws.Range("N136").Value = -12017.8638

# Row 137: H137, J137, L137, N137
$ws.Range("H137").Value = 34088.92
$ws.Range("J137").Value = 34088.92
$ws.Range("L137").Value = 34088.92
$ws.Range("N137").Value = -44288.92

# Row 138: H138, J138, L138, N138
$ws.Range("H138").Value = 25912.062
$ws.Range("J138").Value = 25912.062
$ws.Range("L138").Value = 25912.062
$ws.Range("N138").Value = -36192.06200000001

# Row 139: H139, J139, L139, N139
$ws.Range("H139").Value = 15680.889
$ws.Range("J139").Value = 15680.889
$ws.Range("L139").Value = 15680.889
$ws.Range("N139").Value = -25960.889

# Row 140: H140, J140, L140, N140
$ws.Range("H140").Value = 76554.92999999999
$ws.Range("J140").Value = 76554.92999999999
$ws.Range("L140").Value = 76554.92999999999
$ws.Range("N140").Value = -86914.92999999999

# Row 141: H141, I141, J141, K141, L141, M141, N141
$ws.Range("H141").Value = 34096.332
$ws.Range("I141").Value = 16326
$ws.Range("J141").Value = 42981.5
$ws.Range("K141").Value = 16326
$ws.Range("L141").Value = 42981.5
$ws.Range("M141").Value = -11146
$ws.Range("N141").Value = -53341.5

$ws = $wb.Worksheets.Item("GSM")
# Row 126: H126, I126, J126, K126, L126, M126, N126
$ws.Range("H126").Value = 11112651
$ws.Range("I126").Value = 11112750
$ws.Range("J126").Value = 11112474
$ws.Range("K126").Value = 33338250
$ws.Range("L126").Value = 33337422
$ws.Range("M126").Value = -33335780
$ws.Range("N126").Value = -33342362

# Row 132: H132, I132, J132, K132, L132, M132, N132
$ws.Range("H132").Value = 2632.1025
$ws.Range("I132").Value = 2242.276
$ws.Range("J132").Value = 3762.6
$ws.Range("K132").Value = 6726.828
$ws.Range("L132").Value = 11287.8
$ws.Range("M132").Value = -4196.828
$ws.Range("N132").Value = -16347.8

$ws = $wb.Worksheets.Item("LTW")
# Row 68: H68, I68, K68, M68
$ws.Range("H68").Value = 2857.238
$ws.Range("I68").Value = 2858.9412
$ws.Range("K68").Value = 2858.9412
$ws.Range("M68").Value = -2109.9412

# Row 71: H71, I71, K71, M71
$ws.Range("H71").Value = 2857.238
$ws.Range("I71").Value = 2858.9412
$ws.Range("K71").Value = 14294.706
$ws.Range("M71").Value = -10550.706

$ws = $wb.Worksheets.Item("WVR")
# Row 107: H107, J107, L107, N107
$ws.Range("H107").Value = 494.08334
$ws.Range("J107").Value = 740
$ws.Range("L107").Value = 2220
$ws.Range("N107").Value = -6060
